$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.630.20"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "'3.499.22"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D5").Value = "'598.98"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'180.40"
$ws.Range("E6").Value = "  +4.54%  "
$ws.Range("D7").Value = "'0.615"
$ws.Range("E7").Value = "  +5.46%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'3.500.22"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'0.139"
$ws.Range("E10").Value = "  +5.56%  "
$ws.Range("D11").Value = "'7.02"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("D12").Value = "'0.437"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").Value = "'4.112.28"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'32.36"
$ws.Range("E14").Value = "  +10.52%  "
$ws.Range("D15").Value = "'0.134"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "'67.611.50"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'3.509.22"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'6.33"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'14.33"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "'392.10"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'7.96"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'73.17"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'0.542"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'5.72"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "'0.0000123"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "'10.38"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("E29").Value = "  -2.79%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "'6.16"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "'1.43"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "'2.08"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "'23.60"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'7.44"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'1.62"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "'162.91"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "'0.883"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").Value = "'2.85"
$ws.Range("E40").Value = "  +12.90%  "
$ws.Range("D41").Value = "'1.89"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'6.85"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'26.54"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'2.848.15"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'26.82"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "'0.0725"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").Value = "'41.67"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "'0.0301"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "'335.29"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  -0.55%  "
